$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 (the "Epidemiology Section Editor" row),
# shifting all following rows down by one.
$ws.Rows(3).Insert()

# Fill in the new "Guest Editor" service entry.
$ws.Range("A3").Value = "Guest Editor for Special Issue on Remote Sensing and Crop Health"
$ws.Range("B3").Value = 2019
$ws.Range("C3").Value = "MDPI Remote Sensing"
$ws.Range("D3").Value = "Global"

# The "when" cell of the row pushed down to row 4 keeps the plain/normal
# cell format (rather than the row's bold/emphasis format used by its
# neighbours), matching the author's final formatting.
$ws.Range("A2").Copy()
$ws.Range("B4").PasteSpecial(-4122)

# Match the author's final selection.
$ws.Range("B4").Select()
